$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.1313
$ws.Range("C9").Value = -10.1165
$ws.Range("D12").Value = -7.138999999999999
$ws.Range("D14").Value = -7.687800000000006
$ws.Range("C18").Value = -11.85859999999999
$ws.Range("C20").Value = -11.20390000000001
$ws.Range("D26").Value = -8.6
$ws.Range("C27").Value = -12.59059999999999
$ws.Range("D27").Value = -8.8902
$ws.Range("D29").Value = -7.230200000000001
$ws.Range("C35").Value = -11.24160000000001
$ws.Range("D37").Value = -7.650699999999997
$ws.Range("D38").Value = -8.260799999999994
$ws.Range("D51").Value = -8.110499999999998
$ws.Range("D52").Value = -7.934600000000001
$ws.Range("D55").Value = -8.585799999999999
$ws.Range("C69").Value = -12.2915
$ws.Range("D69").Value = -7.595800000000001
$ws.Range("D70").Value = -7.410600000000001
$ws.Range("C76").Value = -11.98940000000001
$ws.Range("C78").Value = -11.7586
$ws.Range("D81").Value = -7.594500000000002
$ws.Range("C82").Value = -11.4962
$ws.Range("C83").Value = -14.1259
$ws.Range("D83").Value = -8.632499999999999
$ws.Range("C93").Value = -11.0765
$ws.Range("D102").Value = -7.730199999999994
